# ImageSubPortal.xlsx - add an "ISPUrl" column (with hyperlink) in front of
# the existing wrapper-method columns.
#
# The sheet originally has headers in A1:H1 / data in A2:H2. The edit
# inserts a brand-new column A holding:
#   A1 = "ISPUrl"                                             (plain text)
#   A2 = "http://www.esri.com/events/image-submissions#/home" (hyperlink)
# and shifts every existing column one place to the right (B..I).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push all existing columns (A..H) one slot to the right, opening up a
# fresh column A for the new "ISPUrl" field.
$ws.Columns("A").Insert()

# New label cell for the header row.
$ws.Range("A1").Value2 = "ISPUrl"

# New URL cell for the data row, wired up as a real hyperlink (the "#/home"
# fragment becomes the hyperlink's SubAddress/location).
$ws.Hyperlinks.Add($ws.Range("A2"), "http://www.esri.com/events/image-submissions#/home", "/home")

# Match the saved selection state (cursor parked on A3 after the edit).
$ws.Range("A3").Select()
